$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I2 currently holds "TestSuite1()" and I3 holds "TestSuite2()".
# The edit removes the TestSuite1() entry (shifting TestSuite2() up into I2)
# and clears out the now-empty trailing row (I3).
$ws.Range("I2").Value = "TestSuite2()"
$ws.Range("I3").ClearContents()

# Move the active selection to I6 as in the edited workbook.
$ws.Range("I6").Select()
